# Applies the scheduled-runner update to the Leve profit tables (Sheets/Carbuncle_Profits.xlsx).
# For each affected sheet/row, refreshes the market-price / profit columns (H:N) to the new values.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 2201.5
$ws.Range("I6").Value = 15
$ws.Range("J6").Value = 2638.8
$ws.Range("K6").Value = 45
$ws.Range("L6").Value = 7916.400000000001
$ws.Range("M6").Value = 67
$ws.Range("N6").Value = -8140.400000000001
$ws.Range("H7").Value = 30000
$ws.Range("J7").Value = 30000
$ws.Range("L7").Value = 30000
$ws.Range("N7").Value = -30224
$ws.Range("H8").Value = 673
$ws.Range("I8").Value = 673
$ws.Range("K8").Value = 2019
$ws.Range("M8").Value = -1880
$ws.Range("H14").Value = 30000
$ws.Range("J14").Value = 30000
$ws.Range("L14").Value = 30000
$ws.Range("N14").Value = -30382
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").Value = $null
$ws.Range("H38").Value = 660.8570999999999
$ws.Range("I38").Value = 104.333336
$ws.Range("J38").Value = 4000
$ws.Range("K38").Value = 313.000008
$ws.Range("L38").Value = 12000
$ws.Range("M38").Value = 58.99999200000002
$ws.Range("N38").Value = -12744
$ws.Range("H39").Value = 516.1053000000001
$ws.Range("I39").Value = 60
$ws.Range("J39").Value = 679
$ws.Range("K39").Value = 180
$ws.Range("L39").Value = 2037
$ws.Range("M39").Value = 116
$ws.Range("N39").Value = -2629
$ws.Range("H76").Value = 3261.8413
$ws.Range("I76").Value = 3025.8948
$ws.Range("K76").Value = 3025.8948
$ws.Range("M76").Value = -2710.8948
$ws.Range("H79").Value = 3261.8413
$ws.Range("I79").Value = 3025.8948
$ws.Range("K79").Value = 3025.8948
$ws.Range("M79").Value = -1933.8948
$ws.Range("H101").Value = 47619824
$ws.Range("I101").Value = 83333970
$ws.Range("J101").Value = 963.3333
$ws.Range("K101").Value = 250001910
$ws.Range("L101").Value = 2889.9999
$ws.Range("M101").Value = -250000288
$ws.Range("N101").Value = -6133.9999
$ws.Range("H112").Value = 1291.9744
$ws.Range("J112").Value = 1280.7297
$ws.Range("L112").Value = 3842.189100000001
$ws.Range("N112").Value = -6058.189100000001
$ws.Range("H131").Value = 4757.45
$ws.Range("I131").Value = 510
$ws.Range("J131").Value = 4934.4272
$ws.Range("K131").Value = 1530
$ws.Range("L131").Value = 14803.2816
$ws.Range("M131").Value = 3510
$ws.Range("N131").Value = -24883.2816
$ws.Range("H138").Value = 3000.2122
$ws.Range("I138").Value = 1779.3334
$ws.Range("J138").Value = 4017.611
$ws.Range("K138").Value = 5338.0002
$ws.Range("L138").Value = 12052.833
$ws.Range("M138").Value = -198.0002000000004
$ws.Range("N138").Value = -22332.833

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3764.6365
$ws.Range("I32").Value = 1868.403
$ws.Range("K32").Value = 1868.403
$ws.Range("M32").Value = -1581.403
$ws.Range("H45").Value = 2176.4614
$ws.Range("I45").Value = 1499.1428
$ws.Range("J45").Value = 2966.6667
$ws.Range("K45").Value = 1499.1428
$ws.Range("L45").Value = 2966.6667
$ws.Range("M45").Value = -1122.1428
$ws.Range("N45").Value = -3720.6667

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2503.75
$ws.Range("I134").Value = 1472.0667
$ws.Range("J134").Value = 4223.222
$ws.Range("K134").Value = 4416.2001
$ws.Range("L134").Value = 12669.666
$ws.Range("M134").Value = -1881.2001
$ws.Range("N134").Value = -17739.666

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 1750
$ws.Range("I33").Value = 1750
$ws.Range("K33").Value = 1750
$ws.Range("M33").Value = -1371
$ws.Range("H58").Value = 1469.75
$ws.Range("I58").Value = 1650.1666
$ws.Range("J58").Value = 928.5
$ws.Range("K58").Value = 1650.1666
$ws.Range("L58").Value = 928.5
$ws.Range("M58").Value = -1447.1666
$ws.Range("N58").Value = -1334.5
$ws.Range("H68").Value = 8398.75
$ws.Range("J68").Value = 8398.75
$ws.Range("L68").Value = 8398.75
$ws.Range("N68").Value = -9896.75
$ws.Range("H71").Value = 8398.75
$ws.Range("J71").Value = 8398.75
$ws.Range("L71").Value = 25196.25
$ws.Range("N71").Value = -32684.25
$ws.Range("H74").Value = 20314
$ws.Range("J74").Value = 20314
$ws.Range("L74").Value = 20314
$ws.Range("N74").Value = -22062
$ws.Range("H77").Value = 20314
$ws.Range("J77").Value = 20314
$ws.Range("L77").Value = 60942
$ws.Range("N77").Value = -69678
$ws.Range("H136").Value = 1469.75
$ws.Range("I136").Value = 1650.1666
$ws.Range("J136").Value = 928.5
$ws.Range("K136").Value = 4950.4998
$ws.Range("L136").Value = 2785.5
$ws.Range("M136").Value = -2400.4998
$ws.Range("N136").Value = -7885.5

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").Value = $null

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1982.2115
$ws.Range("I122").Value = 1737.381
$ws.Range("J122").Value = 3010.5
$ws.Range("K122").Value = 5212.143
$ws.Range("L122").Value = 9031.5
$ws.Range("M122").Value = -2762.143
$ws.Range("N122").Value = -13931.5

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 13500
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 13500
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 13500
$ws.Range("M32").Value = $null
$ws.Range("N32").Value = -14134
$ws.Range("H93").Value = 9627.182000000001
$ws.Range("I93").Value = 11566
$ws.Range("J93").Value = 902.5
$ws.Range("K93").Value = 11566
$ws.Range("L93").Value = 902.5
$ws.Range("M93").Value = -10318
$ws.Range("N93").Value = -3398.5
$ws.Range("H136").Value = 14495295
$ws.Range("I136").Value = 1272.75
$ws.Range("J136").Value = 22225440
$ws.Range("K136").Value = 3818.25
$ws.Range("L136").Value = 66676320
$ws.Range("M136").Value = -1268.25
$ws.Range("N136").Value = -66681420

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 773.6429000000001
$ws.Range("I81").Value = 628.875
$ws.Range("J81").Value = 966.6667
$ws.Range("K81").Value = 1257.75
$ws.Range("L81").Value = 1933.3334
$ws.Range("M81").Value = -196.75
$ws.Range("N81").Value = -4055.3334
$ws.Range("H84").Value = 773.6429000000001
$ws.Range("I84").Value = 628.875
$ws.Range("J84").Value = 966.6667
$ws.Range("K84").Value = 6288.75
$ws.Range("L84").Value = 9666.666999999999
$ws.Range("M84").Value = -984.75
$ws.Range("N84").Value = -20274.667
$ws.Range("H132").Value = 1839.3784
$ws.Range("I132").Value = 831.7222
$ws.Range("J132").Value = 2794
$ws.Range("K132").Value = 2495.1666
$ws.Range("L132").Value = 8382
$ws.Range("M132").Value = 34.83339999999998
$ws.Range("N132").Value = -13442
